$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep Text type (avoid numeric auto-conversion / trailing-zero loss)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.901.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.300.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.18"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.79%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.31%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.62"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.87"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.77"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.655.38"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.302.85"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.831.65"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.90%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.84"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.06"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.42"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.48"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.98"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.94"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.03"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.94"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.56%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.73"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.012.05"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.72%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.13"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.14"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.28"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.93"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.46"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.523.55"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.60%  "
